$d = $word.ActiveDocument

# 1. Highlight "Technichians by created abrufen" in green
$rng = $d.Content
$rng.Find.Execute("Technichians by created abrufen", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0) | Out-Null
$rng.HighlightColorIndex = 4

# 2. Merge the two runs "Environment abfragen" / " bei Login" (which currently
#    straddle the hidden _GoBack bookmark) into a single run with the new text.
#    Replacing across the bookmark's range removes the old bookmark; it gets
#    re-created later, at the end of the newly appended paragraphs.
$rng2 = $d.Content
$rng2.Find.Execute("Environment abfragen bei Login", $true, $false, $false, $false, $false, `
                    $true, 1, $false, "Environment abfragen bei Login", 2) | Out-Null

# 3. Turn the trailing empty paragraph into "Enable/disable buttons" and append
#    a further paragraph "Beim Login mitarv-Daten mitliefern" after it.
$last = $d.Paragraphs.Item($d.Paragraphs.Count)
$last.Range.InsertBefore("Enable/disable buttons")

$last2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$last2.Range.InsertParagraphAfter()

$last3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$last3.Range.InsertBefore("Beim Login mitarv-Daten mitliefern")

# 4. Re-create the (hidden) _GoBack bookmark as a zero-length range right
#    after the text of the final paragraph. A temporary marker character is
#    used because adding a bookmark exactly at a paragraph's last text
#    position collapses its range to the very start of the document, so the
#    marker keeps the insertion point away from that boundary while the
#    bookmark is created, and is deleted again afterwards.
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalPara.Range.InsertAfter("#")
$markerEnd = $finalPara.Range.End - 1
$bookmarkRange = $d.Range($markerEnd - 1, $markerEnd - 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
$d.Range($markerEnd - 1, $markerEnd).Delete()
